$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure text-like numeric values (e.g. "1.001", "325.83") are written as
# literal text rather than being auto-converted to numbers by Excel, while
# keeping the cell style identical to the original (no style attribute).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.381.71"
Set-TextValue "E2" "  +0.87%  "
Set-TextValue "D3" "1.942.89"
Set-TextValue "E3" "  +2.13%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "325.83"
Set-TextValue "E5" "  +0.38%  "
Set-TextValue "E6" "  +0.15%  "
Set-TextValue "D7" "0.4629"
Set-TextValue "E7" "  +0.71%  "
Set-TextValue "D8" "0.3873"
Set-TextValue "E8" "  -0.33%  "
Set-TextValue "D9" "45.92"
Set-TextValue "E9" "  +0.06%  "
Set-TextValue "D10" "0.07831"
Set-TextValue "E10" "  -0.46%  "
Set-TextValue "D11" "0.9756"
Set-TextValue "E11" "  -1.34%  "
Set-TextValue "D12" "22.64"
Set-TextValue "E12" "  +3.34%  "
Set-TextValue "D13" "1.944.61"
Set-TextValue "E13" "  +3.25%  "
Set-TextValue "D14" "7.084"
Set-TextValue "E14" "  +0.54%  "
Set-TextValue "D15" "5.758"
Set-TextValue "E15" "  -0.23%  "
Set-TextValue "D16" "0.07043"
Set-TextValue "E16" "  +0.43%  "
Set-TextValue "D17" "86.74"
Set-TextValue "E17" "  -1.30%  "
Set-TextValue "D18" "1.003"
Set-TextValue "E18" "  +0.09%  "
Set-TextValue "D19" "0.000009823"
Set-TextValue "E19" "  -0.98%  "
Set-TextValue "E20" "  +0.36%  "
Set-TextValue "D22" "29.409.80"
Set-TextValue "D23" "5.469"
Set-TextValue "E23" "  +2.81%  "
Set-TextValue "E24" "  -0.70%  "
Set-TextValue "D25" "2.162.91"
Set-TextValue "E25" "  +2.29%  "
Set-TextValue "D26" "2.095"
Set-TextValue "E26" "  -0.20%  "
Set-TextValue "D27" "156.94"
Set-TextValue "E27" "  +0.61%  "
Set-TextValue "D28" "19.41"
Set-TextValue "E28" "  -0.21%  "
Set-TextValue "D29" "5.761"
Set-TextValue "E29" "  -2.18%  "
Set-TextValue "D30" "118.36"
Set-TextValue "E30" "  +0.01%  "
Set-TextValue "D31" "1.858"
Set-TextValue "E31" "  -0.67%  "
Set-TextValue "D32" "0.09355"
Set-TextValue "E32" "  +0.20%  "
Set-TextValue "D33" "0.8604"
Set-TextValue "E33" "  -3.68%  "
Set-TextValue "D34" "5.175"
Set-TextValue "E34" "  -1.27%  "
Set-TextValue "E35" "  -1.14%  "
Set-TextValue "E36" "  -0.23%  "
Set-TextValue "D37" "0.05774"
Set-TextValue "E37" "  -0.27%  "
Set-TextValue "E38" "  -0.93%  "
Set-TextValue "D39" "0.02085"
Set-TextValue "E39" "  +0.03%  "
Set-TextValue "D40" "7.689"
Set-TextValue "E40" "  +0.44%  "
Set-TextValue "E41" "  -0.05%  "
Set-TextValue "D42" "0.1780"
Set-TextValue "E42" "  -0.78%  "
Set-TextValue "D43" "9.403"
Set-TextValue "E43" "  -3.10%  "
Set-TextValue "D44" "2.730"
Set-TextValue "E44" "  +7.01%  "
Set-TextValue "D45" "0.000002810"
Set-TextValue "E45" "  +31.92%  "
Set-TextValue "D46" "0.5292"
Set-TextValue "E46" "  -1.08%  "
Set-TextValue "D47" "11.40"
Set-TextValue "E47" "  -3.52%  "
Set-TextValue "B48" "Cronos"
Set-TextValue "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.06866"
Set-TextValue "E48" "  -2.02%  "
Set-TextValue "B49" "RenderToken"
Set-TextValue "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "2.086"
Set-TextValue "E49" "  -5.28%  "
Set-TextValue "E50" "  -1.61%  "
Set-TextValue "D51" "111.34"
Set-TextValue "E51" "  -1.47%  "
